$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (single-dot) need an explicit
# text format first, otherwise Excel auto-converts the string to a Number and the
# canonical OOXML would store it as t="n" instead of text, losing precision/format.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D15", "D17", "D19", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.380.53"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "2.271.53"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "306.57"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "97.07"
$ws.Range("E6").Value = "  +4.87%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "35.70"
$ws.Range("E10").Value = "  +9.65%  "
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D13").Value = "6.66"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "2.630.57"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "14.34"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "2.282.38"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "42.282.94"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "12.47"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "5.95"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "67.61"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "240.47"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "23.80"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "37.59"
$ws.Range("E28").Value = "  +6.28%  "
$ws.Range("D29").Value = "9.49"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "2.10"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "159.02"
$ws.Range("D32").Value = "5.24"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").Value = "0.0739"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").Value = "16.95"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").Value = "0.114"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  +13.89%  "
$ws.Range("D43").Value = "1.998.01"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").Value = "18.80"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "9.96"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").Value = "2.93"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "52.74"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "72.03"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "91.68"
$ws.Range("E51").Value = "  +0.61%  "
